$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: Invalid (G3) and Absent (H3) set to 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Total Attendance Count (D4) and Real (E4) set to 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: Total Attendance Count (D5) and Real (E5) set to 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: Total Attendance Count (D6) and Real (E6) set to 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Rows 7-18: Absent (H) column set to 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
